# Weekly data refresh: insert one new record row for Cilantro at the top
# of the date-ordered table (row 161), pushing every row that currently
# sits at 161-170 down by one (they become 162-171). The dimension and
# formatting expand automatically. The new row copies every column from
# the record that lands below it (the old row 161) and only overrides
# the Fecha (D) and Volumen (J) values with the newly reported figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 161; rows 161..170 shift down to 162..171
$ws.Rows.Item(161).EntireRow.Insert()

# Use the record that is now at row 162 (formerly row 161) as a template
# for the brand-new row 161, so every shared column (Mercado, Region,
# Codreg, Categoria, Variedad, Calidad, precios min/max, unidad, Origen,
# Precio $/Kg, Kg o Unidades, Clasificacion) is populated correctly.
$ws.Range("A162:R162").Copy($ws.Range("A161:R161"))

# Overwrite the two values that differ for this new weekly record
$ws.Cells.Item(161, 4).Value = 44826   # Fecha (D161)
$ws.Cells.Item(161, 10).Value = 3000   # Volumen (J161)
